# "Correcciones de la estructura de la base de datos"
# Adds a new "VIAJE_ASIENTO" table block (rows 151-155) to the data
# dictionary sheet, mirroring the layout of the other table blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed the new block by copying the style skeleton of an existing table
# block (the "Tabla:"/"Objetivo:" pair only spans columns A:B, while the
# header + data rows span the full A:F range) so the new cells pick up
# the exact same cell styles (s="12", s="1", s="5", s="6", s="7", s="8",
# s="9") used throughout the rest of the sheet.
$ws.Range("A139:B140").Copy($ws.Range("A151"))
$ws.Range("A141:F143").Copy($ws.Range("A153"))

# Row 151-152: table name + purpose
$ws.Range("B151").Value = "VIAJE_ASIENTO"
$ws.Range("B152").Value = "Asociar el estado de un asiento en un viaje especifico."

# Row 153 (header) already carries the right values/styles from the copy
# above (Numero / Campo / Tipo de dato / Longitud / Restriccion / Descripcion).

# Row 154: field "asiento"
$ws.Range("A154").Value = 1
$ws.Range("B154").Value = "asiento"
$ws.Range("C154").Value = "Entero"
$ws.Range("D154").Value = ""
$ws.Range("E154").Value = "PK,FK"
$ws.Range("F154").Value = ""

# Row 155: field "viaje"
$ws.Range("A155").Value = 2
$ws.Range("B155").Value = "viaje"
$ws.Range("C155").Value = "Entero"
$ws.Range("D155").Value = ""
$ws.Range("E155").Value = "PK,FK"
$ws.Range("F155").Value = ""

# Reflect the author's final cursor position in the sheet view.
$ws.Range("C103").Select()
